# Refactor weight handling methods in DenseLayer and NeuralNetwork
#
# The re-trained model produces new predictions, which ripple through
# this sliding-window results sheet:
#   C (IPC PO)   = new model prediction for each point
#   D (DELTA)    = C - B            (prediction error)
#   E (DELTA^2)  = D^2              (squared error)
#   Row 52 TOTAL = SUM(D2:D51) in C52, SUM(E2:E51) in E52
#   Row 53 MSE   = E52 / 50         (mean squared error)
#
# Values below are written directly (matching the exact recomputed
# model output) rather than re-derived in-sheet, so results are exact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (C: IPC PO, D: DELTA, E: DELTA^2)
$rowData = @{
    2 = @(29.34289985487003, -0.2071001451299708, 0.04289047011285495)
    3 = @(29.73399442021041, -0.01600557978958861, 0.0002561785844008874)
    4 = @(30.30976685091175, 0.4697668509117463, 0.2206808942155389)
    5 = @(30.18441060497628, 0.3744106049762763, 0.1401833011187012)
    6 = @(29.91739787915503, -0.002602120844972688, 0.000006771032891841378196075372)
    7 = @(29.86322763495295, -0.1167723650470549, 0.01363578523868265)
    8 = @(29.87031407170072, -0.1696859282992769, 0.02879331426278733)
    9 = @(29.98186659486311, -0.2281334051368944, 0.05204485053935439)
    10 = @(30.13287623834016, -0.08712376165983926, 0.007590549845760478)
    11 = @(30.19635820737606, -0.1836417926239413, 0.03372430799813467)
    12 = @(30.49092004476196, 0.05092004476195555, 0.002592850958559557)
    13 = @(30.41667902625066, -0.06332097374933809, 0.004009545716564364)
    14 = @(30.41725733225557, -0.2727426677444349, 0.07438856280835121)
    15 = @(30.75600763974486, 0.006007639744858295, 0.000036091735304001040709337422)
    16 = @(30.8328102252049, -0.1071897747950992, 0.01148964782062408)
    17 = @(30.92742909302136, -0.02257090697863973, 0.0005094458418384075)
    18 = @(31.13534542000329, 0.115345420003294, 0.01330456591573628)
    19 = @(31.24816772355857, 0.1281677235585654, 0.01642696536218485)
    20 = @(31.1295982748876, -0.1504017251124026, 0.02262067891678671)
    21 = @(30.98313353305045, -0.3968664669495467, 0.1575029925890157)
    22 = @(31.2016208500709, -0.3783791499290956, 0.143170781101065)
    23 = @(31.66159758224181, 0.01159758224181573, 0.0001345039138556796)
    24 = @(32.49035952587511, 0.6103595258751078, 0.3725387508264864)
    25 = @(32.47948710295107, 0.1994871029510676, 0.03979510424380986)
    26 = @(32.84956581219205, 0.3995658121920442, 0.159652838272688)
    27 = @(33.13122302259087, 0.2812230225908721, 0.07908638843514613)
    28 = @(33.19406796210896, 0.2940679621089615, 0.0864759663389176)
    29 = @(33.31829812935527, 0.2182981293552686, 0.04765407328000958)
    30 = @(33.71936165410794, 0.3193616541079436, 0.1019918661145618)
    31 = @(33.72889839495205, 0.02889839495205138, 0.0008351172308047485)
    32 = @(33.79738590116477, -0.3026140988352282, 0.09157529281385723)
    33 = @(34.20263769160752, -0.1973623083924778, 0.03895188077400752)
    34 = @(34.45974062913849, -0.4402593708615115, 0.1938283136313739)
    35 = @(35.62825551613716, 0.3282555161371619, 0.1077516838744745)
    36 = @(35.89664932593882, 0.1966493259388145, 0.03867095739219008)
    37 = @(36.17513560476671, -0.1248643952332884, 0.01559111719697486)
    38 = @(36.79703670043323, -0.002963299566765443, 0.000008781144322392261098819514)
    39 = @(36.88889165364485, -0.4111083463551495, 0.1690100724428655)
    40 = @(37.99193202877537, 0.09193202877536777, 0.008451497914755048)
    41 = @(38.69646754624609, 0.1964675462460903, 0.03859949672795963)
    42 = @(39.00851112573505, 0.1085111257350562, 0.01177466440828917)
    43 = @(39.37586542242229, -0.02413457757771198, 0.0005824778348545979)
    44 = @(39.6478100977316, -0.2521899022683982, 0.06359974680614425)
    45 = @(39.85032822565149, -0.2496717743485135, 0.06233599490633505)
    46 = @(40.03887235974106, -0.5611276402589382, 0.3148642286625644)
    47 = @(40.34758395275934, -0.552416047240655, 0.3051634892489896)
    48 = @(41.43800678165266, 0.2380067816526577, 0.05664722811265586)
    49 = @(41.53548273587013, 0.03548273587012574, 0.001259024544829108)
    50 = @(41.95881490504632, 0.1588149050463272, 0.02522217406487393)
    51 = @(42.78848311819058, 0.5884831181905739, 0.3463123803953009)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
}

# TOTAL row (52): sum of DELTA (C52) and sum of DELTA^2 (E52)
$ws.Cells.Item(52, 3).Value = -0.07116797080472992
$ws.Cells.Item(52, 5).Value = 3.764223663269034

# MSE row (53): mean of DELTA^2 over the 50 points
$ws.Cells.Item(53, 5).Value = 0.07528447326538068
